$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.207, 0.372, 0.245, -19, 80, 79),
    @(0.112, -0.118, 0.473, 74, 33, 166),
    @(-0.001, 0.224, 0.5610000000000001, -56, -60, 62),
    @(-0.174, -0.015, 0.283, -20, -12, 98),
    @(-0.081, -0.034, 0.11, -71, 7, -150),
    @(0.003, 0.101, 0.265, -168, 76, -4),
    @(0.212, 0.089, 0.197, -11, -20, -110),
    @(0.02, 0.032, 0.263, 107, -80, -70),
    @(0.042, 0.177, 0.089, -75, -70, -43),
    @(0.04, 0.144, 0.284, 172, -5, -30)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    for ($c = 0; $c -lt $values[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}
